# Applies the cryptos.xlsx price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.998.47'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.868.44'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.31%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08172'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.19'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('E11').Value = '  -2.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.71'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.21%  '
$ws.Range('D13').Value = '1.870.39'
$ws.Range('E13').Value = '  -2.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.246'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.145'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.81%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.61'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.90%  '
$ws.Range('E18').Value = '  -3.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06378'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.07%  '
$ws.Range('E20').Value = '  -1.47%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '29.973.58'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.788'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.46%  '
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.162'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.80%  '
$ws.Range('D26').Value = '2.083.14'
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.220'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.89'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.61%  '
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1033'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.900'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.734'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02410'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.213'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06327'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2135'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.171'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.483'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6284'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.203'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.90'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5877'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.625'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.202'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.140'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.44%  '
